# Scheduled runner update: refresh market-price derived figures (currentAveragePrice*,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) for the affected leve rows across each crafting sheet.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ALC_updates = @(
    @{ Addr = "H6"; Val = 39.444443 },
    @{ Addr = "I6"; Val = 45 },
    @{ Addr = "J6"; Val = 20 },
    @{ Addr = "K6"; Val = 135 },
    @{ Addr = "L6"; Val = 60 },
    @{ Addr = "M6"; Val = -23 },
    @{ Addr = "N6"; Val = -284 },
    @{ Addr = "H28"; Val = 2096.625 },
    @{ Addr = "I28"; Val = 907.2778 },
    @{ Addr = "K28"; Val = 907.2778 },
    @{ Addr = "M28"; Val = -422.2778 },
    @{ Addr = "H53"; Val = 1085.6316 },
    @{ Addr = "I53"; Val = 1664.8334 },
    @{ Addr = "J53"; Val = 818.3077 },
    @{ Addr = "K53"; Val = 1664.8334 },
    @{ Addr = "L53"; Val = 818.3077 },
    @{ Addr = "M53"; Val = -1027.8334 },
    @{ Addr = "N53"; Val = -2092.3077 },
    @{ Addr = "H58"; Val = 4577.625 },
    @{ Addr = "I58"; Val = 520.8 },
    @{ Addr = "K58"; Val = 1562.4 },
    @{ Addr = "M58"; Val = -1412.4 },
    @{ Addr = "H134"; Val = 72342.62 },
    @{ Addr = "J134"; Val = 72342.62 },
    @{ Addr = "L134"; Val = 72342.62 },
    @{ Addr = "N134"; Val = -82482.62 },
    @{ Addr = "H136"; Val = 78212 },
    @{ Addr = "J136"; Val = 78212 },
    @{ Addr = "L136"; Val = 78212 },
    @{ Addr = "N136"; Val = -88412 }
)
foreach ($u in $ALC_updates) {
    $ws.Range($u.Addr).Value = $u.Val
}

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ARM_updates = @(
    @{ Addr = "H7"; Val = 26713.166 },
    @{ Addr = "J7"; Val = 26713.166 },
    @{ Addr = "L7"; Val = 26713.166 },
    @{ Addr = "N7"; Val = -26941.166 },
    @{ Addr = "H52"; Val = 53479.75 },
    @{ Addr = "J52"; Val = 53479.75 },
    @{ Addr = "L52"; Val = 53479.75 },
    @{ Addr = "N52"; Val = -54115.75 },
    @{ Addr = "H74"; Val = 36601.207 },
    @{ Addr = "I74"; Val = 42893.125 },
    @{ Addr = "J74"; Val = 6400 },
    @{ Addr = "K74"; Val = 42893.125 },
    @{ Addr = "L74"; Val = 6400 },
    @{ Addr = "M74"; Val = -42019.125 },
    @{ Addr = "N74"; Val = -8148 },
    @{ Addr = "H77"; Val = 36601.207 },
    @{ Addr = "I77"; Val = 42893.125 },
    @{ Addr = "J77"; Val = 6400 },
    @{ Addr = "K77"; Val = 214465.625 },
    @{ Addr = "L77"; Val = 32000 },
    @{ Addr = "M77"; Val = -210097.625 },
    @{ Addr = "N77"; Val = -40736 },
    @{ Addr = "H102"; Val = 159015.58 },
    @{ Addr = "I102"; Val = 251649.75 },
    @{ Addr = "J102"; Val = 35503.332 },
    @{ Addr = "K102"; Val = 251649.75 },
    @{ Addr = "L102"; Val = 35503.332 },
    @{ Addr = "M102"; Val = -250027.75 },
    @{ Addr = "N102"; Val = -38747.332 },
    @{ Addr = "H118"; Val = 49220 },
    @{ Addr = "J118"; Val = 49220 },
    @{ Addr = "L118"; Val = 49220 },
    @{ Addr = "N118"; Val = -52534 },
    @{ Addr = "H121"; Val = 54282.855 },
    @{ Addr = "J121"; Val = 54282.855 },
    @{ Addr = "L121"; Val = 54282.855 },
    @{ Addr = "N121"; Val = -57776.855 },
    @{ Addr = "H132"; Val = 1935.1923 },
    @{ Addr = "I132"; Val = 1316.4 },
    @{ Addr = "J132"; Val = 3997.8333 },
    @{ Addr = "K132"; Val = 3949.2 },
    @{ Addr = "L132"; Val = 11993.4999 },
    @{ Addr = "M132"; Val = -1419.2 },
    @{ Addr = "N132"; Val = -17053.4999 }
)
foreach ($u in $ARM_updates) {
    $ws.Range($u.Addr).Value = $u.Val
}

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$BSM_updates = @(
    @{ Addr = "H52"; Val = 99990 },
    @{ Addr = "J52"; Val = 99990 },
    @{ Addr = "L52"; Val = 99990 },
    @{ Addr = "N52"; Val = -100516 },
    @{ Addr = "H86"; Val = 1280.5834 },
    @{ Addr = "J86"; Val = 1014.5714 },
    @{ Addr = "L86"; Val = 1014.5714 },
    @{ Addr = "N86"; Val = -3260.5714 },
    @{ Addr = "H89"; Val = 1280.5834 },
    @{ Addr = "J89"; Val = 1014.5714 },
    @{ Addr = "L89"; Val = 5072.857 },
    @{ Addr = "N89"; Val = -16304.857 },
    @{ Addr = "H110"; Val = 67396.60000000001 },
    @{ Addr = "J110"; Val = 67396.60000000001 },
    @{ Addr = "L110"; Val = 67396.60000000001 },
    @{ Addr = "N110"; Val = -75576.60000000001 },
    @{ Addr = "H114"; Val = 89989.2 },
    @{ Addr = "J114"; Val = 89989.2 },
    @{ Addr = "L114"; Val = 89989.2 },
    @{ Addr = "N114"; Val = -98667.2 },
    @{ Addr = "H115"; Val = 85278.71000000001 },
    @{ Addr = "J115"; Val = 94990 },
    @{ Addr = "L115"; Val = 94990 },
    @{ Addr = "N115"; Val = -98124 },
    @{ Addr = "H116"; Val = 83742.5 },
    @{ Addr = "J116"; Val = 83742.5 },
    @{ Addr = "L116"; Val = 83742.5 },
    @{ Addr = "N116"; Val = -92920.5 },
    @{ Addr = "H117"; Val = 96167.71000000001 },
    @{ Addr = "J117"; Val = 96167.71000000001 },
    @{ Addr = "L117"; Val = 96167.71000000001 },
    @{ Addr = "N117"; Val = -105345.71 },
    @{ Addr = "H121"; Val = 99990 },
    @{ Addr = "J121"; Val = 99990 },
    @{ Addr = "L121"; Val = 99990 },
    @{ Addr = "N121"; Val = -103484 },
    @{ Addr = "H122"; Val = 80716.25 },
    @{ Addr = "J122"; Val = 80716.25 },
    @{ Addr = "L122"; Val = 80716.25 },
    @{ Addr = "N122"; Val = -90516.25 },
    @{ Addr = "H127"; Val = 54983.875 },
    @{ Addr = "J127"; Val = 54983.875 },
    @{ Addr = "L127"; Val = 54983.875 },
    @{ Addr = "N127"; Val = -64903.875 },
    @{ Addr = "H129"; Val = 80000 },
    @{ Addr = "J129"; Val = 80000 },
    @{ Addr = "L129"; Val = 80000 },
    @{ Addr = "N129"; Val = -90000 }
)
foreach ($u in $BSM_updates) {
    $ws.Range($u.Addr).Value = $u.Val
}

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$CRP_updates = @(
    @{ Addr = "H58"; Val = 1821.5 },
    @{ Addr = "I58"; Val = 1737.5834 },
    @{ Addr = "J58"; Val = 2325 },
    @{ Addr = "K58"; Val = 1737.5834 },
    @{ Addr = "L58"; Val = 2325 },
    @{ Addr = "M58"; Val = -1534.5834 },
    @{ Addr = "N58"; Val = -2731 },
    @{ Addr = "H60"; Val = 12812.692 },
    @{ Addr = "I60"; Val = 6415 },
    @{ Addr = "K60"; Val = 6415 },
    @{ Addr = "M60"; Val = -5904 },
    @{ Addr = "H114"; Val = 48745 },
    @{ Addr = "J114"; Val = 48745 },
    @{ Addr = "L114"; Val = 48745 },
    @{ Addr = "N114"; Val = -57423 },
    @{ Addr = "H118"; Val = 61452.11 },
    @{ Addr = "J118"; Val = 61452.11 },
    @{ Addr = "L118"; Val = 61452.11 },
    @{ Addr = "N118"; Val = -64766.11 },
    @{ Addr = "H136"; Val = 1821.5 },
    @{ Addr = "I136"; Val = 1737.5834 },
    @{ Addr = "J136"; Val = 2325 },
    @{ Addr = "K136"; Val = 5212.7502 },
    @{ Addr = "L136"; Val = 6975 },
    @{ Addr = "M136"; Val = -2662.7502 },
    @{ Addr = "N136"; Val = -12075 }
)
foreach ($u in $CRP_updates) {
    $ws.Range($u.Addr).Value = $u.Val
}

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$GSM_updates = @(
    @{ Addr = "H108"; Val = 54410.54 },
    @{ Addr = "J108"; Val = 54410.54 },
    @{ Addr = "L108"; Val = 54410.54 },
    @{ Addr = "N108"; Val = -62090.54 },
    @{ Addr = "H109"; Val = 80265 },
    @{ Addr = "J109"; Val = 80265 },
    @{ Addr = "L109"; Val = 80265 },
    @{ Addr = "N109"; Val = -82345 },
    @{ Addr = "H110"; Val = 80655.62 },
    @{ Addr = "J110"; Val = 80655.62 },
    @{ Addr = "L110"; Val = 80655.62 },
    @{ Addr = "N110"; Val = -88835.62 },
    @{ Addr = "H117"; Val = 85739.60000000001 },
    @{ Addr = "J117"; Val = 85739.60000000001 },
    @{ Addr = "L117"; Val = 85739.60000000001 },
    @{ Addr = "N117"; Val = -92623.60000000001 },
    @{ Addr = "H118"; Val = 47950 },
    @{ Addr = "I118"; Val = 0 },
    @{ Addr = "J118"; Val = 47950 },
    @{ Addr = "K118"; Val = 0 },
    @{ Addr = "L118"; Val = 47950 },
    @{ Addr = "N118"; Val = -51264 }
)
foreach ($u in $GSM_updates) {
    $ws.Range($u.Addr).Value = $u.Val
}
$ws.Range("M118").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$LTW_updates = @(
    @{ Addr = "H16"; Val = 660.7179599999999 },
    @{ Addr = "I16"; Val = 623.1818 },
    @{ Addr = "K16"; Val = 623.1818 },
    @{ Addr = "M16"; Val = -453.1818 },
    @{ Addr = "H118"; Val = 56780.4 },
    @{ Addr = "J118"; Val = 56780.4 },
    @{ Addr = "L118"; Val = 56780.4 },
    @{ Addr = "N118"; Val = -60094.4 },
    @{ Addr = "H131"; Val = 82141.22 },
    @{ Addr = "I131"; Val = 51765.332 },
    @{ Addr = "J131"; Val = 97329.164 },
    @{ Addr = "K131"; Val = 51765.332 },
    @{ Addr = "L131"; Val = 97329.164 },
    @{ Addr = "M131"; Val = -46725.332 },
    @{ Addr = "N131"; Val = -107409.164 },
    @{ Addr = "H132"; Val = 3035.1667 },
    @{ Addr = "J132"; Val = 3468.3333 },
    @{ Addr = "L132"; Val = 10404.9999 },
    @{ Addr = "N132"; Val = -15464.9999 }
)
foreach ($u in $LTW_updates) {
    $ws.Range($u.Addr).Value = $u.Val
}

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$WVR_updates = @(
    @{ Addr = "H113"; Val = 3624.5 },
    @{ Addr = "I113"; Val = 3624.5 },
    @{ Addr = "K113"; Val = 10873.5 },
    @{ Addr = "M113"; Val = -8703.5 },
    @{ Addr = "H121"; Val = 44886.668 },
    @{ Addr = "J121"; Val = 44886.668 },
    @{ Addr = "L121"; Val = 44886.668 },
    @{ Addr = "N121"; Val = -48380.668 },
    @{ Addr = "H132"; Val = 1013046.8 },
    @{ Addr = "I132"; Val = 1839.5769 },
    @{ Addr = "J132"; Val = 2559599 },
    @{ Addr = "K132"; Val = 5518.7307 },
    @{ Addr = "L132"; Val = 7678797 },
    @{ Addr = "M132"; Val = -2988.7307 },
    @{ Addr = "N132"; Val = -7683857 },
    @{ Addr = "H136"; Val = 2503.25 },
    @{ Addr = "J136"; Val = 2726.3125 },
    @{ Addr = "L136"; Val = 8178.9375 },
    @{ Addr = "N136"; Val = -13278.9375 }
)
foreach ($u in $WVR_updates) {
    $ws.Range($u.Addr).Value = $u.Val
}
